$d = $word.ActiveDocument

# ---- Pass 1: replace each source text with a unique temporary marker (breaks the move-cycles) ----
$d.Content.Find.Execute("Oferecimento de seminários aos alunos sobre temas atuais de  Ciências, Engenharia e Empreendedorismo.", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_0_0zZ", 2) | Out-Null
$d.Content.Find.Execute("Seminários abrangendo os cenários atuais e futuros da indústria de alta tecnologia e do campo de atuação do engenheiro físico.", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_0_1zZ", 2) | Out-Null
$d.Content.Find.Execute("Os seminários proferidos por estudantes de graduação e pós-graduação, professores e convidados serão debatidos e analisados pelos alunos em forma de relatório. Os seminários apresentados pelos alunos serão avaliados na disciplina.", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_0_2zZ", 2) | Out-Null
$d.Content.Find.Execute("Não há.", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_0_3zZ", 2) | Out-Null
$d.Content.Find.Execute("3577649 - Carlos Angelo Nunes", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_0_4zZ", 2) | Out-Null
$d.Content.Find.Execute("1176388 - Luiz Tadeu Fernandes Eleno", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_1_0zZ", 2) | Out-Null
$d.Content.Find.Execute("Seminários seguido de debates com profissionais e estudantes de graduação e pós-graduação sobre temas relevantes e atuais das áreas de Física, Tecnologia e Engenharia, abrangendo desde as pesquisas básicas até o segmento industrial e de serviços.", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_1_1zZ", 2) | Out-Null
$d.Content.Find.Execute("A nota final será calculada pela média aritmética dos relatórios e do seminário.", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_1_2zZ", 2) | Out-Null
$d.Content.Find.Execute("A ser definido de acordo com os temas dos seminários.", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_1_3zZ", 2) | Out-Null
$d.Content.Find.Execute("Offering seminars to students on current science, engineering and entrepreneurship topics.", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_2_0zZ", 2) | Out-Null
$d.Content.Find.Execute("Seminars covering the current and future scenarios of the high technology industry and the field of activity of the physical engineer.", $true, $false, $false, $false, $false, $true, 1, $false, "ZzMarkerLOM3243_2_1zZ", 2) | Out-Null

# ---- Pass 2: replace each marker with its final destination text ----
$d.Content.Find.Execute("ZzMarkerLOM3243_0_0zZ", $true, $false, $false, $false, $false, $true, 1, $false, "Seminários abrangendo os cenários atuais e futuros da indústria de alta tecnologia e do campo de atuação do engenheiro físico.", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_0_1zZ", $true, $false, $false, $false, $false, $true, 1, $false, "Os seminários proferidos por estudantes de graduação e pós-graduação, professores e convidados serão debatidos e analisados pelos alunos em forma de relatório. Os seminários apresentados pelos alunos serão avaliados na disciplina.", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_0_2zZ", $true, $false, $false, $false, $false, $true, 1, $false, "Não há.", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_0_3zZ", $true, $false, $false, $false, $false, $true, 1, $false, "3577649 - Carlos Angelo Nunes", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_0_4zZ", $true, $false, $false, $false, $false, $true, 1, $false, "Oferecimento de seminários aos alunos sobre temas atuais de  Ciências, Engenharia e Empreendedorismo.", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_1_0zZ", $true, $false, $false, $false, $false, $true, 1, $false, "Seminários seguido de debates com profissionais e estudantes de graduação e pós-graduação sobre temas relevantes e atuais das áreas de Física, Tecnologia e Engenharia, abrangendo desde as pesquisas básicas até o segmento industrial e de serviços.", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_1_1zZ", $true, $false, $false, $false, $false, $true, 1, $false, "A nota final será calculada pela média aritmética dos relatórios e do seminário.", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_1_2zZ", $true, $false, $false, $false, $false, $true, 1, $false, "A ser definido de acordo com os temas dos seminários.", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_1_3zZ", $true, $false, $false, $false, $false, $true, 1, $false, "1176388 - Luiz Tadeu Fernandes Eleno", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_2_0zZ", $true, $false, $false, $false, $false, $true, 1, $false, "Seminars covering the current and future scenarios of the high technology industry and the field of activity of the physical engineer.", 2) | Out-Null
$d.Content.Find.Execute("ZzMarkerLOM3243_2_1zZ", $true, $false, $false, $false, $false, $true, 1, $false, "Offering seminars to students on current science, engineering and entrepreneurship topics.", 2) | Out-Null
